$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in previously-empty monthly figures for Jul-Nov 2023
$ws.Range("H2").Value = 0.6589
$ws.Range("I2").Value = 0.7171
$ws.Range("J2").Value = 0.6136
$ws.Range("K2").Value = 0.6061
$ws.Range("L2").Value = 0.5779

# Row 2: set the source-link cell (previously blank styled cell) and drop its style
$ws.Range("O2").ClearFormats()
$ws.Range("O2").Value = "http://www.yahii.com.br/poupanca.html"

# Rows 3-5: shift the source-link values down a row; the oldest one falls off
$ws.Range("O3").Value = "http://www.acinh.com.br/servicos/indicadores-economicos/poupanca-mensal"
$ws.Range("O4").Value = "https://www.portalbrasil.net/poupanca_mensal/"
$ws.Range("O5").ClearContents()
